# Rename the "Input [Source Name]" and "Output [Source Name]" table/header
# columns to "Input [Sample Name]" / "Output [Sample Name]" on the
# Events-Tillage sheet's annotationTable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events-Tillage")
$tbl = $ws.ListObjects.Item("annotationTable")

# Column 1: "Input [Source Name]" -> "Input [Sample Name]"
$tbl.ListColumns.Item(1).Range.Rows.Item(1).Value = "Input [Sample Name]"

# Column 18: "Output [Source Name]" -> "Output [Sample Name]"
$tbl.ListColumns.Item(18).Range.Rows.Item(1).Value = "Output [Sample Name]"
